# Update the cryptos list snapshot (GitHub Actions scheduled refresh).
# Applies the latest Price (column D) and Volume(1h) (column E) figures,
# plus the HuobiToken / InjectiveProtocol row swap (rows 49-50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "41.783.62"
$ws.Range("E2").Value = "  +5.40%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.229.72"
$ws.Range("E3").Value = "  +2.75%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.34"
$ws.Range("E5").Value = "  +2.08%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.01%  "

# Row 7 - Solana
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.88"
$ws.Range("E7").Value = "  -2.11%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.03%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +2.66%  "

# Row 10 - OKB
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.34"
$ws.Range("E10").Value = "  +1.57%  "

# Row 11 - Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0893"
$ws.Range("E11").Value = "  +4.95%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.30%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.562.36"
$ws.Range("E13").Value = "  +2.88%  "

# Row 14 - Chainlink
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.63"
$ws.Range("E14").Value = "  -1.78%  "

# Row 15 - Avalanche
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.97"
$ws.Range("E15").Value = "  +0.80%  "

# Row 16 - Polygon
$ws.Range("E16").Value = "  -1.14%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  +1.77%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.250.63"
$ws.Range("E18").Value = "  +3.83%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "41.759.27"
$ws.Range("E19").Value = "  +5.49%  "

# Row 20 - Litecoin
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.20"
$ws.Range("E20").Value = "  +0.67%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  -2.39%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.05"
$ws.Range("E22").Value = "  +0.70%  "

# Row 23 - BitcoinCash
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.66"
$ws.Range("E23").Value = "  +8.74%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.03%  "

# Row 25
$ws.Range("E25").Value = "  +2.33%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.37"
$ws.Range("E26").Value = "  +1.94%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.69"
$ws.Range("E27").Value = "  +2.10%  "

# Row 28 - Kaspa
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.143"
$ws.Range("E28").Value = "  +2.76%  "

# Row 29 - Monero
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.89"
$ws.Range("E29").Value = "  -2.36%  "

# Row 30 - EthereumClassic
$ws.Range("E30").Value = "  +0.58%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  -2.42%  "

# Row 33 - Stellar
$ws.Range("E33").Value = "  -0.55%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.99"
$ws.Range("E34").Value = "  +6.31%  "

# Row 35 - Filecoin
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.68"
$ws.Range("E35").Value = "  +3.53%  "

# Row 36 - Hedera
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0635"
$ws.Range("E36").Value = "  +2.92%  "

# Row 37 - THORChain
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.64"
$ws.Range("E37").Value = "  -4.73%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  -4.67%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  -1.26%  "

# Row 40 - TerraClassic
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.000255"
$ws.Range("E40").Value = "  +30.04%  "

# Row 41 - BinanceUSD
$ws.Range("E41").Value = "  -0.03%  "

# Row 42 - VeChain
$ws.Range("E42").Value = "  +4.48%  "

# Row 43 - FTXToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.86"
$ws.Range("E43").Value = "  -0.87%  "

# Row 44 - FraxShare
$ws.Range("E44").Value = "  +8.44%  "

# Row 45 - Cronos
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0979"
$ws.Range("E45").Value = "  +6.22%  "

# Row 46 - TrustWalletToken
$ws.Range("E46").Value = "  +0.36%  "

# Row 47 - Aave
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.91"
$ws.Range("E47").Value = "  -3.74%  "

# Row 48 - Maker
$ws.Range("D48").Value = "1.480.24"
$ws.Range("E48").Value = "  -2.15%  "

# Row 49 - was InjectiveProtocol, now HuobiToken (row swap with 50)
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.81"
$ws.Range("E49").Value = "  +0.23%  "

# Row 50 - was HuobiToken, now InjectiveProtocol (row swap with 49)
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.46"
$ws.Range("E50").Value = "  -7.03%  "

# Row 51 - MultiversX
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.30"
$ws.Range("E51").Value = "  +5.46%  "
